# Auto-generated edit script: updates Coin/Link/Price/Volume(1h) cells per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "257.50"
Set-TextCell "E2" "4.96%"
Set-TextCell "D3" "27.57"
Set-TextCell "E3" "-2.90%"
Set-TextCell "D4" "5.219"
Set-TextCell "E4" "-0.56%"
Set-TextCell "D5" "0.05924"
Set-TextCell "E5" "3.87%"
Set-TextCell "E6" "0.90%"
Set-TextCell "D7" "0.8677"
Set-TextCell "D8" "1.028"
Set-TextCell "E8" "13.56%"
Set-TextCell "B9" "One"
Set-TextCell "C9" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D9" "0.01051"
Set-TextCell "E9" "1,668.34%"
Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1420"
Set-TextCell "E10" "3.67%"
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.07195"
Set-TextCell "E11" "1.84%"
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.03264"
Set-TextCell "E12" "3.46%"
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.09224"
Set-TextCell "E13" "-0.01%"
Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001541"
Set-TextCell "E14" "0.36%"
Set-TextCell "D15" "0.005935"
Set-TextCell "E15" "-0.21%"
Set-TextCell "D16" "3.486"
Set-TextCell "E16" "-0.10%"
Set-TextCell "D17" "3.270"
Set-TextCell "E17" "2.13%"
Set-TextCell "D18" "2.206"
Set-TextCell "E18" "1.41%"
Set-TextCell "E19" "-0.64%"
Set-TextCell "D20" "0.03605"
Set-TextCell "E20" "9.46%"
Set-TextCell "E21" "1.94%"
Set-TextCell "D22" "3.525"
Set-TextCell "E22" "0.00%"
Set-TextCell "D23" "0.04161"
Set-TextCell "E23" "2.14%"
Set-TextCell "E24" "1.57%"
Set-TextCell "D25" "0.001218"
Set-TextCell "E25" "-0.29%"
Set-TextCell "D26" "0.004525"
Set-TextCell "E26" "8.92%"
Set-TextCell "E27" "0.10%"
Set-TextCell "E28" "33.91%"
Set-TextCell "D40" "0.03823"
Set-TextCell "E40" "1.21%"
Set-TextCell "D41" "0.005413"
Set-TextCell "E41" "4.25%"
Set-TextCell "D42" "0.1104"
Set-TextCell "E42" "3.70%"
Set-TextCell "D43" "0.001901"
Set-TextCell "E43" "-13.54%"
Set-TextCell "D44" "0.009827"
Set-TextCell "E44" "7.47%"
Set-TextCell "D45" "0.00005434"
Set-TextCell "E45" "3.19%"
Set-TextCell "E47" "4.02%"
Set-TextCell "D48" "0.002162"
Set-TextCell "E48" "-4.69%"
